$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue 'D2' '67.332.22'
$ws.Range('E2').Value = '  -3.25%  '
Set-TextValue 'D3' '3.495.21'
$ws.Range('E3').Value = '  -4.82%  '
$ws.Range('E4').Value = '  +0.01%  '
Set-TextValue 'D5' '606.15'
$ws.Range('E5').Value = '  -2.47%  '
Set-TextValue 'D6' '148.74'
$ws.Range('E6').Value = '  -6.55%  '
Set-TextValue 'D7' '3.493.28'
$ws.Range('E7').Value = '  -4.76%  '
$ws.Range('E8').Value = '  -0.02%  '
Set-TextValue 'D9' '0.481'
$ws.Range('E9').Value = '  -3.18%  '
$ws.Range('E10').Value = '  -3.84%  '
Set-TextValue 'D11' '6.96'
$ws.Range('E11').Value = '  -2.82%  '
$ws.Range('E12').Value = '  -3.95%  '
$ws.Range('E13').Value = '  -4.64%  '
Set-TextValue 'D14' '4.084.24'
$ws.Range('E14').Value = '  -4.74%  '
Set-TextValue 'D15' '31.43'
$ws.Range('E15').Value = '  -2.86%  '
Set-TextValue 'D16' '3.499.18'
$ws.Range('E16').Value = '  -4.26%  '
Set-TextValue 'D17' '67.200.76'
$ws.Range('E17').Value = '  -3.48%  '
$ws.Range('E18').Value = '  -0.62%  '
$ws.Range('E19').Value = '  -1.90%  '
Set-TextValue 'D20' '15.01'
$ws.Range('E20').Value = '  -5.60%  '
Set-TextValue 'D21' '446.11'
$ws.Range('E21').Value = '  -5.15%  '
Set-TextValue 'D22' '9.02'
$ws.Range('E22').Value = '  -12.59%  '
Set-TextValue 'D23' '0.621'
$ws.Range('E23').Value = '  -4.28%  '
Set-TextValue 'D24' '77.15'
$ws.Range('E24').Value = '  -3.31%  '
Set-TextValue 'D25' '0.0000129'
$ws.Range('E25').Value = '  +4.65%  '
$ws.Range('E26').Value = '  +0.09%  '
Set-TextValue 'D27' '3.632.49'
$ws.Range('E27').Value = '  -4.84%  '
Set-TextValue 'D28' '10.14'
$ws.Range('E28').Value = '  -8.16%  '
Set-TextValue 'D29' '8.28'
Set-TextValue 'D30' '2.48'
$ws.Range('E30').Value = '  -4.69%  '
$ws.Range('E32').Value = '  +0.05%  '
Set-TextValue 'D33' '0.162'
$ws.Range('E33').Value = '  -0.39%  '
Set-TextValue 'D34' '25.67'
$ws.Range('E34').Value = '  -3.38%  '
$ws.Range('E35').Value = '  -4.05%  '
$ws.Range('E36').Value = '  -6.48%  '
Set-TextValue 'D37' '3.482.59'
$ws.Range('E37').Value = '  -5.12%  '
$ws.Range('E39').Value = '  +0.08%  '
Set-TextValue 'D40' '0.999'
$ws.Range('E40').Value = '  +0.02%  '
Set-TextValue 'D41' '173.50'
$ws.Range('E41').Value = '  -2.41%  '
$ws.Range('E42').Value = '  -1.32%  '
Set-TextValue 'D43' '0.0875'
$ws.Range('E43').Value = '  -1.75%  '
Set-TextValue 'D44' '5.42'
$ws.Range('E44').Value = '  -6.35%  '
Set-TextValue 'D45' '0.881'
$ws.Range('E45').Value = '  -4.69%  '
Set-TextValue 'D46' '45.44'
$ws.Range('E46').Value = '  -2.62%  '
$ws.Range('B47').Value = 'ONDO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue 'D47' '1.27'
$ws.Range('E47').Value = '  +4.99%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D48' '27.04'
$ws.Range('E48').Value = '  -6.29%  '
Set-TextValue 'D49' '2.57'
$ws.Range('E49').Value = '  -5.28%  '
$ws.Range('E50').Value = '  -4.10%  '
$ws.Range('E51').Value = '  -3.32%  '
